$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. "Object Code" sheet: append four new rows (Index 11, 12, 13, 20)
# ------------------------------------------------------------------
$objectCode = $wb.Worksheets.Item("Object Code")

$objectCode.Range("A7").Value = 11

$objectCode.Range("A8").Value = 12
$objectCode.Range("B8").Value = "add wall"

$objectCode.Range("A9").Value = 13
$objectCode.Range("B9").Value = "remove wall"

$objectCode.Range("A10").Value = 20
$objectCode.Range("B10").Value = "interpolating entity"

$objectCode.Range("A11").Value = 21

# ------------------------------------------------------------------
# 2. Insert three new worksheets between "Object Code" and "Image ID":
#      Client Player Message, Other Player Message, Bullet Message
# ------------------------------------------------------------------
$clientMsg = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $objectCode)
$clientMsg.Name = "Client Player Message"

$otherMsg = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $clientMsg)
$otherMsg.Name = "Other Player Message"

$bulletMsg = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $otherMsg)
$bulletMsg.Name = "Bullet Message"

# Re-fetch "Image ID" only now (after the inserts above) so the
# reference is not stale for the formatting Copy() calls below.
$imageId = $wb.Worksheets.Item("Image ID")

# ------------------------------------------------------------------
# Helper formatting source: "Image ID" B1 carries the bold/underline
# header style (cell style index 6) that all three new header rows use
# on *both* columns.
# ------------------------------------------------------------------

# ---------------- Client Player Message ----------------
$imageId.Range("B1").Copy()
$clientMsg.Range("A1:B1").PasteSpecial(-4122)

$clientMsg.Range("A1").Value = "Index"
$clientMsg.Range("B1").Value = "Data"

$clientMsg.Range("A2").Value = 0
$clientMsg.Range("B2").Value = '"000"'
$clientMsg.Range("A3").Value = 1
$clientMsg.Range("B3").Value = "highest snapshot"
$clientMsg.Range("A4").Value = 2
$clientMsg.Range("B4").Value = "client id"
$clientMsg.Range("A5").Value = 3
$clientMsg.Range("B5").Value = "image id"
$clientMsg.Range("A6").Value = 4
$clientMsg.Range("B6").Value = "sprite index"
$clientMsg.Range("A7").Value = 5
$clientMsg.Range("B7").Value = "x"
$clientMsg.Range("A8").Value = 6
$clientMsg.Range("B8").Value = "y"
$clientMsg.Range("A9").Value = 7
$clientMsg.Range("B9").Value = "w"
$clientMsg.Range("A10").Value = 8
$clientMsg.Range("B10").Value = "h"
$clientMsg.Range("A11").Value = 9
$clientMsg.Range("B11").Value = "rotation"
$clientMsg.Range("A12").Value = 10
$clientMsg.Range("B12").Value = "alpha"
$clientMsg.Range("A13").Value = 11
$clientMsg.Range("B13").Value = "role"
$clientMsg.Range("A14").Value = 12
$clientMsg.Range("B14").Value = "team"
$clientMsg.Range("A15").Value = 13
$clientMsg.Range("B15").Value = "current weapon"
$clientMsg.Range("A16").Value = 14
$clientMsg.Range("B16").Value = "item(empty/imageid)"
$clientMsg.Range("A17").Value = 15
$clientMsg.Range("B17").Value = "health"
$clientMsg.Range("A18").Value = 16
$clientMsg.Range("B18").Value = "is_invincible"
$clientMsg.Range("A19").Value = 17
$clientMsg.Range("B19").Value = "speed_boost"
$clientMsg.Range("A20").Value = 18
$clientMsg.Range("B20").Value = "damage_boost"
$clientMsg.Range("A21").Value = 19
$clientMsg.Range("B21").Value = "visibility"

$clientMsg.Columns.Item(2).ColumnWidth = 19.8
$clientMsg.Range("A1:B12").Select()

# ---------------- Other Player Message ----------------
$imageId.Range("B1").Copy()
$otherMsg.Range("A1:B1").PasteSpecial(-4122)

$otherMsg.Range("A1").Value = "Index"
$otherMsg.Range("B1").Value = "Data"

$otherMsg.Range("A2").Value = 0
$otherMsg.Range("B2").Value = '"000"'
$otherMsg.Range("A3").Value = 1
$otherMsg.Range("B3").Value = "client id"
$otherMsg.Range("A4").Value = 2
$otherMsg.Range("B4").Value = "image id"
$otherMsg.Range("A5").Value = 3
$otherMsg.Range("B5").Value = "sprite index"
$otherMsg.Range("A6").Value = 4
$otherMsg.Range("B6").Value = "x"
$otherMsg.Range("A7").Value = 5
$otherMsg.Range("B7").Value = "y"
$otherMsg.Range("A8").Value = 6
$otherMsg.Range("B8").Value = "w"
$otherMsg.Range("A9").Value = 7
$otherMsg.Range("B9").Value = "h"
$otherMsg.Range("A10").Value = 8
$otherMsg.Range("B10").Value = "rotation"
$otherMsg.Range("A11").Value = 9
$otherMsg.Range("B11").Value = "alpha"

$otherMsg.Columns.Item(2).ColumnWidth = 17.8
$otherMsg.Range("B3").Select()

# ---------------- Bullet Message ----------------
$imageId.Range("B1").Copy()
$bulletMsg.Range("A1:B1").PasteSpecial(-4122)

$bulletMsg.Range("A1").Value = "Index"
$bulletMsg.Range("B1").Value = "Data"

$bulletMsg.Range("A2").Value = 0
$bulletMsg.Range("B2").Value = '"001"'
$bulletMsg.Range("A3").Value = 1
$bulletMsg.Range("B3").Value = "client id"
$bulletMsg.Range("A4").Value = 2
$bulletMsg.Range("B4").Value = "image id"
$bulletMsg.Range("A5").Value = 3
$bulletMsg.Range("B5").Value = "sprite index"
$bulletMsg.Range("A6").Value = 4
$bulletMsg.Range("B6").Value = "x"
$bulletMsg.Range("A7").Value = 5
$bulletMsg.Range("B7").Value = "y"
$bulletMsg.Range("A8").Value = 6
$bulletMsg.Range("B8").Value = "w"
$bulletMsg.Range("A9").Value = 7
$bulletMsg.Range("B9").Value = "h"
$bulletMsg.Range("A10").Value = 8
$bulletMsg.Range("B10").Value = "rotation"
$bulletMsg.Range("A11").Value = 9
$bulletMsg.Range("B11").Value = "alpha"
$bulletMsg.Range("A12").Value = 10
$bulletMsg.Range("B12").Value = "speed"
$bulletMsg.Range("A13").Value = 11
$bulletMsg.Range("B13").Value = "damage"
$bulletMsg.Range("A14").Value = 12
$bulletMsg.Range("B14").Value = "team"
$bulletMsg.Range("A15").Value = 13
$bulletMsg.Range("B15").Value = "xRatio"
$bulletMsg.Range("A16").Value = 14
$bulletMsg.Range("B16").Value = "yRatio"

$bulletMsg.Columns.Item(2).ColumnWidth = 18.6
$bulletMsg.Range("B14").Select()

# ------------------------------------------------------------------
# 3. Restore "Object Code" as the active sheet/selection
# ------------------------------------------------------------------
$objectCode.Activate()
$objectCode.Range("B10").Select()
